$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1. Add the three new character styles referenced by the diff.
#    They must be created as character styles (type 2) BEFORE we
#    assign them to any Range.Style, so the runtime correctly
#    serializes them as <w:rStyle> inside the run's <w:rPr> instead
#    of as a paragraph style.
# ---------------------------------------------------------------

$gaNStyle = $d.Styles.Add("GaNStyle", 2)
$gaNStyle.Font.NameAscii = "Calibri"
$gaNStyle.Font.Name = "Calibri"
$gaNStyle.Font.Size = 14

$gaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$gaNParagraph.Font.NameAscii = "Calibri"
$gaNParagraph.Font.Name = "Calibri"
$gaNParagraph.Font.Size = 10

$gaNLinks = $d.Styles.Add("GaNLinks", 2)
$gaNLinks.Font.NameAscii = "Calibri"
$gaNLinks.Font.Name = "Calibri"
$gaNLinks.Font.Bold = $true
$gaNLinks.Font.Color = 8388608
$gaNLinks.Font.Size = 9.5
$gaNLinks.Font.Underline = 1

# ---------------------------------------------------------------
# 2. Fix the duplicated trailing date text in the four campaign
#    paragraphs ("...Herkules souhvezdi.<dup dates>" -> "...Herkules souhvezdi.")
# ---------------------------------------------------------------

$oldTail = "Herkules souhvězdí.13. – 22. června, 12. – 21. července, 10. – 19. srpna"
$newTail = "Herkules souhvězdí."
$fixRange = $d.Content
$fixRange.Find.ClearFormatting()
$fixRange.Find.Execute($oldTail, $true, $false, $false, $false, $false, $true, 1, $false, $newTail, 2)

# ---------------------------------------------------------------
# 3. Apply the GaNStyle character style to the run containing the
#    whole campaign paragraph text (all four occurrences).
# ---------------------------------------------------------------

$target = "Informace v této příručce jsou určeny pro pozorovací kampaň probíhající od 13. – 22. června, 12. – 21. července, 10. – 19. srpna. Při pozorování použijte hvězdy oblohy, které zobrazujíHerkules souhvězdí."

$searchRange = $d.Content
$searchRange.Find.ClearFormatting()
$searchRange.Find.Text = $target
$searchRange.Find.Forward = $true
$searchRange.Find.Wrap = 1

while ($searchRange.Find.Execute()) {
    $searchRange.Style = "GaNStyle"
    $searchRange.Collapse(0)
}

# ---------------------------------------------------------------
# 4. Apply the GaNLinks character style to the "Jenik Hollan" credit
#    run.
# ---------------------------------------------------------------

$linkTarget = "Jeník Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/"
$linkRange = $d.Content
$linkRange.Find.ClearFormatting()
$linkRange.Find.Text = $linkTarget
$linkRange.Find.Forward = $true
$linkRange.Find.Wrap = 1

while ($linkRange.Find.Execute()) {
    $linkRange.Style = "GaNLinks"
    $linkRange.Collapse(0)
}

Write-Output "done"
